$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to append (social, prescription, script #, fill date, exp. Date, refills avail)
$rows = @(
    @("812-121-0912", "Lozol",      9382398232,  43104, 43469, 3),
    @("812-121-0912", "Prozac",     11212435324, 42883, 43982, 18),
    @("503-388-1908", "Kezvara",    4739518341,  43004, 43391, 0),
    @("503-388-1908", "Diuril",     47394710247, 43364, 43729, 11),
    @("503-388-1908", "Anafril",    430983092,   41339, 43449, 2),
    @("810-225-7205", "Zaraoxolyn", 28737124212, 43202, 43567, 9)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Apply a date number format to the fill date / exp. date columns, reusing a
# single style definition for every cell (mirrors the workbook being edited
# once in Excel and the format being applied to the whole block at once).
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Copy()
$ws.Range("D2:E7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Columns C, D and E are sized to best-fit their (now populated) contents.
$ws.Columns.Item(3).ColumnWidth = 11.166666666666666
$ws.Columns.Item(4).ColumnWidth = 8.721354166666666
$ws.Columns.Item(5).ColumnWidth = 9.721354166666666

$ws.Range("D6").Select()
